$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# Helper: replace the run content of a whole paragraph (leaving the
# paragraph's own <w:pPr> / paragraph-mark untouched) with a hand built
# sequence of <w:r> elements supplied as raw OOXML. This lets us control
# run boundaries exactly, the same way Word's editor would leave them
# after a human typed/edited the text in several passes.
# ----------------------------------------------------------------------
function Set-ParagraphRuns($para, [string]$runsXml) {
    $start = $para.Range.Start
    $end = $para.Range.End - 1   # exclude the trailing paragraph mark
    $sub = $d.Range($start, $end)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
           $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $sub.InsertXML($xml)
}

function Find-ParagraphByText([string]$text) {
    foreach ($para in $d.Paragraphs) {
        $t = $para.Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $para
        }
    }
    return $null
}

# ------------------------------------------------------------------
# 1) Rewrite the "Created a Sign-Up, Sign-In, and Sign-Out ..." bullet
#    with the updated wording from the commit, split across runs the
#    same way the author's edit left them.
# ------------------------------------------------------------------
$oldBullet = "Created a Sign-Up, Sign-In, and Sign-Out handler using typescript so that the users will be able to get, post, put and delete their data in the database."
$bulletPara = Find-ParagraphByText $oldBullet

if ($bulletPara -ne $null) {
    $rPr = '<w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>'
    $runsXml =
        '<w:r>' + $rPr + '<w:t>Created a Sign-Up, Sign-In, Sign-Out</w:t></w:r>' +
        '<w:r>' + $rPr + '<w:t xml:space="preserve"> handler </w:t></w:r>' +
        '<w:r>' + $rPr + '<w:t xml:space="preserve">using typescript so that end users can </w:t></w:r>' +
        '<w:r>' + $rPr + '<w:t>in</w:t></w:r>' +
        '<w:r>' + $rPr + '<w:t xml:space="preserve">sert their </w:t></w:r>' +
        '<w:r>' + $rPr + '<w:t>data in the database.</w:t></w:r>'
    Set-ParagraphRuns $bulletPara $runsXml
}

# ------------------------------------------------------------------
# 2) Move the (rendering-only) lastRenderedPageBreak marker off the
#    "Implemented Contact Me page ..." run and onto the following
#    "Deep Dive Coding" heading run, matching the re-paginated layout
#    that results from the bullet above getting shorter.
# ------------------------------------------------------------------
$contactPara = Find-ParagraphByText "Implemented Contact Me page so that end user can contact me."
if ($contactPara -ne $null) {
    $rPr = '<w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>'
    $runsXml = '<w:r>' + $rPr + '<w:t>Implemented Contact Me page so that end user can contact me.</w:t></w:r>'
    Set-ParagraphRuns $contactPara $runsXml
}

$deepDivePara = Find-ParagraphByText "Deep Dive Coding"
if ($deepDivePara -ne $null) {
    $rPr = '<w:rPr><w:b/><w:color w:val="2F5496" w:themeColor="accent1" w:themeShade="BF"/><w:spacing w:val="20"/></w:rPr>'
    $runsXml = '<w:r>' + $rPr + '<w:lastRenderedPageBreak/><w:t>Deep Dive Coding</w:t></w:r>'
    Set-ParagraphRuns $deepDivePara $runsXml
}

Write-Host "Edit complete"
